# Coastal Surface Piercing Profilers - update ref des from CP05MOAS-GL005 to
# CP05MOAS-GL388 on both the Moorings and Asset_Cal_Info sheets, and make
# Asset_Cal_Info the active sheet/tab.

$wb = $excel.ActiveWorkbook

$wsMoorings = $wb.Worksheets.Item("Moorings")
$wsCalInfo  = $wb.Worksheets.Item("Asset_Cal_Info")

# --- Moorings sheet: Ref Des for the deployment record ---
$wsMoorings.Range("A2").Value = "CP05MOAS-GL388"

# --- Asset_Cal_Info sheet: Ref Des for each instrument's calibration block ---
$wsCalInfo.Range("A2").Value  = "CP05MOAS-GL388-01-ADCPAM000"
$wsCalInfo.Range("A3").Value  = "CP05MOAS-GL388-01-ADCPAM000"
$wsCalInfo.Range("A4").Value  = "CP05MOAS-GL388-01-ADCPAM000"
$wsCalInfo.Range("A5").Value  = "CP05MOAS-GL388-01-ADCPAM000"

$wsCalInfo.Range("A7").Value  = "CP05MOAS-GL388-02-FLORTM000"
$wsCalInfo.Range("A8").Value  = "CP05MOAS-GL388-02-FLORTM000"
$wsCalInfo.Range("A9").Value  = "CP05MOAS-GL388-02-FLORTM000"
$wsCalInfo.Range("A10").Value = "CP05MOAS-GL388-02-FLORTM000"

$wsCalInfo.Range("A12").Value = "CP05MOAS-GL388-03-CTDGVM000"
$wsCalInfo.Range("A14").Value = "CP05MOAS-GL388-04-DOSTAM000"
$wsCalInfo.Range("A16").Value = "CP05MOAS-GL388-05-PARADM000"
$wsCalInfo.Range("A18").Value = "CP05MOAS-GL388-00-ENG000000"

# --- Make Asset_Cal_Info the active/selected tab (was Moorings) ---
$wsCalInfo.Activate()
